$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Done" status for the two newly completed functions (rows 5 and 6)
$ws.Range("F5").Value = "DOne"
$ws.Range("F6").Value = "Done"

# Update the active selection to reflect the last edited cell
$ws.Range("F7").Select()
